$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 35159.619052
$ws.Range("B4").Value = 1340.040996
$ws.Range("B9").Value = 0.995319
$ws.Range("B10").Value = 44.70195
$ws.Range("B12").Value = 44171.326773
$ws.Range("B14").Value = 6.260751
$ws.Range("C14").Value = 13071.753763
$ws.Range("B15").Value = 2199.204126
$ws.Range("B16").Value = 251.181764
$ws.Range("B17").Value = 75313.301416
$ws.Range("C17").Value = 298329.594594
$ws.Range("B18").Value = 203.008032
$ws.Range("B20").Value = 6829.221888
$ws.Range("B25").Value = 1388.752039
$ws.Range("B27").Value = 644.472393
$ws.Range("C27").Value = 104307.02145
$ws.Range("B29").Value = 1216.879861
$ws.Range("B30").Value = 298.546878
$ws.Range("B35").Value = 80362.961505
$ws.Range("B37").Value = 8744.111525
$ws.Range("B43").Value = 9865.732379
$ws.Range("B44").Value = 2669.999126
$ws.Range("B46").Value = 30123.551377
$ws.Range("B49").Value = 4634.382571
$ws.Range("B51").Value = 0.535716
$ws.Range("C51").Value = 57716
$ws.Range("B58").Value = 170.997826
$ws.Range("B65").Value = 53.038586
$ws.Range("C65").Value = 12419.533725
$ws.Range("B70").Value = 371.148249
$ws.Range("B77").Value = 4.75
$ws.Range("B81").Value = 34.248201
$ws.Range("B84").Value = 579.35
$ws.Range("B90").Value = 6.4
$ws.Range("B91").Value = 21.719858
